# Updates odds/score values on Sheet1 to reflect the latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "G2"  = 2.9
    "L2"  = 4
    "AH2" = 12
    "AQ2" = 67

    "N4" = 8
    "Q4" = 2.3
    "R4" = 1.6

    "G7"  = 2.05
    "I7"  = 3.4
    "J7"  = 2.63
    "AC7" = 15
    "AH7" = 19
    "AJ7" = 34
    "AK7" = 23
    "AM7" = 126
    "AX7" = 17
    "BB7" = 126

    "O8" = 1.2
    "P8" = 4.33

    "G15"  = 3.1
    "I15"  = 2.55
    "K15"  = 1.83
    "L15"  = 3.4
    "M15"  = 1.13
    "N15"  = 6
    "Y15"  = 13
    "AA15" = 34
    "AG15" = 6
    "AR15" = 126
    "AX15" = 15
    "BB15" = 351

    "G16"  = 2
    "H16"  = 3
    "I16"  = 4.33
    "L16"  = 5
    "U16"  = 2.25
    "V16"  = 1.57
    "X16"  = 8
    "Z16"  = 17
    "AE16" = 21
    "AG16" = 9
    "AX16" = 26
    "AZ16" = 101

    "H17"  = 3.2
    "K17"  = 2
    "N17"  = 9
    "Q17"  = 2.2
    "R17"  = 1.65
    "AA17" = 26
    "AC17" = 8
    "AH17" = 12
    "AM17" = 351

    "O18"  = 1.57
    "P18"  = 2.25
    "U18"  = 2.25
    "V18"  = 1.57
    "AE18" = 21
    "AG18" = 6.5
    "BA18" = 126
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
